$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Preserve the original column A values (A1:A45) before we touch rows,
    # because this edit only removes data from columns B:I of row 3
    # (shifting the rows below it up) while column A's literal index values
    # stay exactly as they were for rows 1-45; only row 46 disappears entirely.
    $lastRow = 46
    $savedA = @()
    for ($r = 1; $r -le ($lastRow - 1); $r++) {
        $savedA += , $ws.Cells.Item($r, 1).Value()
    }

    # Delete worksheet row 3 entirely, shifting rows 4-46 up to 3-45.
    $ws.Rows.Item(3).Delete()

    # Restore column A (1) to its original, un-shifted values for rows 1-45.
    for ($r = 1; $r -le ($lastRow - 1); $r++) {
        $ws.Cells.Item($r, 1).Value = $savedA[$r - 1]
    }
}
